$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "EKD - Coalition Party  (Eestil Koonderakond, EKD)"
$ws.Range("C1").Value = "ERL - People's Union  of Estonia (Eestimaa Rahvaliit, ERL)"
$ws.Range("D1").Value = "EÜRP - United People's Party of Estonia (Eestimaa Ühendatud Rahvapartei, EÜRP)"
$ws.Range("E1").Value = "IL - Pro Patria Union (Isamaaliit, IL)"
$ws.Range("F1").Value = "KeE - Centre Party (Keskerakond, KeE)"
$ws.Range("G1").Value = "M - Moderates  (Mõõdukad, M)"
$ws.Range("H1").Value = "RE - Reform Party   (, RE)"
$ws.Range("I1").Value = "RL - People's Union  (, RL)"
$ws.Range("J1").Value = "RP - Res Publica (Res Publica, RP)"
$ws.Range("K1").Value = "ER - Estonian Greens (Eestimaa Rohelised, ER)"
$ws.Range("L1").Value = "IRL - Pro Patria and Res Publica Union (Isamaa ja Res Publica Liit, IRL)"
$ws.Range("M1").Value = "SDTP - Social Democratic Labour Party (, SDTP)"
$ws.Range("N1").Value = "EKRE - Conservative People’s Party (Eesti Konservatiivne Rahvaerakond, EKRE)"
$ws.Range("O1").Value = "EVA - Estonian Free Party (Eesti Vabaerakond, EVA)"
$ws.Range("P1").Value = "SDE - Social Democratic Party (, SDE)"
